$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("OBSERVATION_FAC")
$ws4 = $wb.Worksheets.Item("BIRD_SPECIE_DIM")

Write-Host "C2 before:" $ws1.Range("C2").Value
Write-Host "F2 before:" $ws1.Range("F2").Value
Write-Host "I3 before:" $ws1.Range("I3").Value

$ws1.Range("I3").Value = "DEFAULT FALSE"
$ws1.Range("C2").Value = "FOREIGN KEY"
$ws1.Range("F2").Value = "FOREIGN KEY"

Write-Host "C2 after:" $ws1.Range("C2").Value
Write-Host "F2 after:" $ws1.Range("F2").Value
Write-Host "I3 after:" $ws1.Range("I3").Value

$ws1.Activate()
$ws1.Range("E4").Select()

$ws4.Activate()
Write-Host "activesheet" $wb.ActiveSheet.Name
